$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 222 - this shifts rows 222:298 down to 223:299
# (matches dimension growing from A1:R298 to A1:R299)
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new weekly record.
# Values mirror the former row 222 except for Fecha (D) and Volumen (J),
# which carry the new week's data.
$ws.Cells.Item(222, 1).Value = 7
$ws.Cells.Item(222, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(222, 3).Value = "Ñuble"
$ws.Cells.Item(222, 4).Value = 45146
$ws.Cells.Item(222, 5).Value = 16
$ws.Cells.Item(222, 6).Value = 100112040
$ws.Cells.Item(222, 7).Value = "Cilantro"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 180
$ws.Cells.Item(222, 11).Value = 1500
$ws.Cells.Item(222, 12).Value = 1500
$ws.Cells.Item(222, 13).Value = 1500
$ws.Cells.Item(222, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(222, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(222, 16).Value = 1500
$ws.Cells.Item(222, 17).Value = 1
$ws.Cells.Item(222, 18).Value = "Hortaliza"
